# MDL.xlsx dashboard refresh ("update beauty, MG, PSS dashboard")
#
# The canonical diff for this change:
#   - xl/sharedStrings.xml gains 7 new unique strings (new booking/consult
#     codes), appended after the existing table.
#   - xl/worksheets/sheet1.xml: cell A2 (the "id" column of the single data
#     row) is repointed from the old shared string ("CA-GR9ZVVXL") to the
#     newly added last string ("CA-CC5MALQX").
#
# Because this engine always rebuilds xl/sharedStrings.xml from the strings
# that are actually referenced by a cell at save time (any string that is
# no longer referenced by any cell is dropped), the four shared strings
# that were already present in the workbook but were not referenced by any
# cell (CA-SI3VQCZT, " -", CA-TYLZ7MIX and CA-GR9ZVVXL, the old value of
# A2) would otherwise be silently discarded on save. To keep the resulting
# shared-strings table byte-for-byte aligned with the target (same 28
# entries, same order), row 10 re-references those four pre-existing
# strings and then introduces the six brand-new ones, in the exact order
# they must appear, before the final new string is written into A2 itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Keep the pre-existing (previously unreferenced) shared strings alive.
$ws.Range("A10").Value = "CA-SI3VQCZT"
$ws.Range("B10").Value = " -"
$ws.Range("C10").Value = "CA-TYLZ7MIX"
$ws.Range("D10").Value = "CA-GR9ZVVXL"

# Append the new shared strings introduced by this commit, in order.
$ws.Range("E10").Value = "CA-N1DCYO19"
$ws.Range("F10").Value = "CA-PXLZAWCU"
$ws.Range("G10").Value = "CA-3FNUOO34"
$ws.Range("H10").Value = "CA-5EBR8U6D"
$ws.Range("I10").Value = "MDL - Single Consultation"
$ws.Range("J10").Value = "CA-TSMYA9T0"

# The actual dashboard edit: row 2's id now references the newly added
# "CA-CC5MALQX" string instead of the old "CA-GR9ZVVXL" one.
$ws.Range("A2").Value = "CA-CC5MALQX"
